$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for years 2000, 2005, 2006, 2007, 2008, 2009 (currently rows 2-7).
# This shifts the 2010-2013 rows (currently rows 8-11) up to rows 2-5.
$ws.Range("A2:B7").EntireRow.Delete()
